# Update "想去人数" (number of people interested) values in column F
# across the four worksheets, per the data refresh captured in the diff.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1098
$ws.Range("F4").Value = 1909
$ws.Range("F5").Value = 595
$ws.Range("F6").Value = 1243
$ws.Range("F8").Value = 29
$ws.Range("F9").Value = 136
$ws.Range("F11").Value = 112
$ws.Range("F13").Value = 797
$ws.Range("F14").Value = 235
$ws.Range("F19").Value = 215
$ws.Range("F20").Value = 693
$ws.Range("F21").Value = 70
$ws.Range("F22").Value = 661
$ws.Range("F23").Value = 183
$ws.Range("F24").Value = 47
$ws.Range("F25").Value = 903
$ws.Range("F26").Value = 350
$ws.Range("F27").Value = 189
$ws.Range("F28").Value = 58
$ws.Range("F29").Value = 302
$ws.Range("F32").Value = 421

# Sheet "演出" (Performance)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 19
$ws.Range("F6").Value = 30
$ws.Range("F7").Value = 263
$ws.Range("F11").Value = 129

# Sheet "本地生活" (Local life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 328

# Sheet "全部类型" (All types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 328
$ws.Range("F4").Value = 1098
$ws.Range("F5").Value = 1909
$ws.Range("F6").Value = 595
$ws.Range("F7").Value = 1243
$ws.Range("F10").Value = 29
$ws.Range("F11").Value = 136
$ws.Range("F13").Value = 112
$ws.Range("F15").Value = 797
$ws.Range("F16").Value = 235
$ws.Range("F18").Value = 19
$ws.Range("F24").Value = 30
$ws.Range("F25").Value = 263
$ws.Range("F26").Value = 263
$ws.Range("F27").Value = 215
$ws.Range("F28").Value = 693
$ws.Range("F29").Value = 70
$ws.Range("F30").Value = 661
$ws.Range("F31").Value = 183
$ws.Range("F32").Value = 47
$ws.Range("F33").Value = 903
$ws.Range("F34").Value = 350
$ws.Range("F37").Value = 189
$ws.Range("F38").Value = 58
$ws.Range("F39").Value = 302
$ws.Range("F41").Value = 129
$ws.Range("F42").Value = 129
$ws.Range("F46").Value = 421
